# Feat/49/transfer crud
# Adds a new "TRANSFER" worksheet between ACCOUNT and REGULAR_TRANSFER, and
# adds a TRANSFER_ID column to TEMPORARY_TRANSFER.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) TEMPORARY_TRANSFER: add a new "TRANSFER_ID" column (G) tying each
#    temporary transfer row back to a transfer record.
#    (Done first so the "TRANSFER_ID" shared string is registered before the
#    strings used on the new TRANSFER sheet.)
# ---------------------------------------------------------------------------
$tempWs = $wb.Worksheets.Item("TEMPORARY_TRANSFER")

$tempWs.Range("G1").Value = "TRANSFER_ID"
$tempWs.Range("G2").Value = 1
$tempWs.Range("G3").Value = 2
$tempWs.Range("G4").Value = 3

$tempWs.Rows.Item(1).RowHeight = 28
$tempWs.Rows.Item(2).RowHeight = 14
$tempWs.Rows.Item(3).RowHeight = 14
$tempWs.Rows.Item(4).RowHeight = 14

$tempWs.Range("G53").Select()

# ---------------------------------------------------------------------------
# 2) Insert the new TRANSFER worksheet before REGULAR_TRANSFER.
# ---------------------------------------------------------------------------
$regularWs = $wb.Worksheets.Item("REGULAR_TRANSFER")
$transferWs = $wb.Worksheets.Add($regularWs)
$transferWs.Name = "TRANSFER"

# Header row
$transferWs.Range("A1").Value = "ID"
$transferWs.Range("B1").NumberFormat = "@"
$transferWs.Range("B1").Value = "TITLE"
$transferWs.Range("C1").Value = "USER_ID"

# Row 2
$transferWs.Range("A2").Value = 1
$transferWs.Range("B2").NumberFormat = "@"
$transferWs.Range("B2").Value = "2023/06"
$transferWs.Range("C2").Value = 1

# Row 3
$transferWs.Range("A3").Value = 2
$transferWs.Range("B3").NumberFormat = "@"
$transferWs.Range("B3").Value = "2023/07"
$transferWs.Range("C3").Value = 1

# Row 4
$transferWs.Range("A4").Value = 3
$transferWs.Range("B4").NumberFormat = "@"
$transferWs.Range("B4").Value = "2022/09"
$transferWs.Range("C4").Value = 2

$transferWs.Rows.Item(1).RowHeight = 14
$transferWs.Rows.Item(2).RowHeight = 14
$transferWs.Rows.Item(3).RowHeight = 14
$transferWs.Rows.Item(4).RowHeight = 14

$transferWs.Range("F16").Select()
